$d = $word.ActiveDocument

# Simple text corrections (typo fixes "gerência" -> "gerencia", and one semantic change)
$d.Content.Find.Execute("Gerente gerência representantes", $true, $false, $false, $false, $false, $true, 1, $false, "Gerente gerencia representantes", 2)
$d.Content.Find.Execute("Gerente gerência clientes", $true, $false, $false, $false, $false, $true, 1, $false, "Gerente gerencia clientes", 2)
$d.Content.Find.Execute("Gerente gerência template dos contratos de serviços", $true, $false, $false, $false, $false, $true, 1, $false, "Gerente gerencia template dos contratos de serviços", 2)
$d.Content.Find.Execute("Gerente acompanha contratos", $true, $false, $false, $false, $false, $true, 1, $false, "Gerente gerencia contratos", 2)
$d.Content.Find.Execute("Gerente gerência serviços", $true, $false, $false, $false, $false, $true, 1, $false, "Gerente gerencia serviços", 2)

# RF8 row description text change
$d.Content.Find.Execute("Representante gerência relatórios de serviço", $true, $false, $false, $false, $false, $true, 1, $false, "Gerente visualiza solicitações de serviço do cliente", 2)

